# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the
# e674504d-fc0a-4e2d-a9f8-9da8f52936e6.md row after re-generating the
# handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G4").Value = "2016-08-18 22:47:58"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-18 22:47:53"
$wsZhCn.Range("K4").Value = "2016-08-18 22:48:17"

# de-de sheet: Correspond Handoff Datetime (shared with Overview's value) /
# Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-08-18 22:47:58"
$wsDeDe.Range("K4").Value = "2016-08-18 22:48:24"
